# Fruta / hortaliza, semanal
# New weekly price observation is inserted at row 20, pushing the existing
# rows 20-47 down to 21-48 (and the sheet dimension grows to R48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 20; this shifts rows 20:47 -> 21:48
# and carries their formatting with them, exactly like a manual Excel
# "Insert Sheet Rows" on row 20.
$ws.Rows("20").Insert()

# Populate the freshly-inserted row 20 with the new observation.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 45070
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112043
$ws.Range("G20").Value = "Pepino dulce"
$ws.Range("H20").Value = "Cultivar IV Región"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14600
$ws.Range("N20").Value = "$/bandeja 18 kilos"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 811
$ws.Range("Q20").Value = 18
$ws.Range("R20").Value = "Hortaliza"
